$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Replace the numeric value in A2 with the text "MODELO"
$ws.Range("A2").Value = "MODELO"

# Move the active selection to C6, matching the saved cursor position
$ws.Range("C6").Select()
